# Adds two new trailing columns (U: DF_chg, V: HO_chg) to the Indonesia
# shapefile-adding worksheet, mirroring the IF_val/IF_chg (S/T) columns that
# already exist: most provinces get a 0 placeholder styled like column R/S/T
# (style index 18, NumberFormat "0.00"), while the handful of provinces that
# already carry real FTT/IF numbers also get real DF_chg / HO_chg numbers
# (those particular cells are left in the default/general style, matching
# how columns S/T/T themselves are formatted for those same rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new column titles ------------------------------
$ws.Range("U1").Value = "DF_chg"
$ws.Range("V1").Value = "HO_chg"

# Cells in U2:U35 / V2:V35 that must end up with the "0.00" numeric style
# (index 18) and a literal 0 -- i.e. every row except the ones with real
# FTT-derived data below (and, for U, except U10/U14/U27 which carry real
# numbers themselves).
$uStyledZeroRows = 2..35 | Where-Object { $_ -ne 10 -and $_ -ne 14 -and $_ -ne 27 }
foreach ($r in $uStyledZeroRows) {
    $cell = $ws.Range("U$r")
    $cell.Value = 0
    $cell.NumberFormat = "0.00"
}

$vStyledZeroRows = 2..35 | Where-Object { $_ -ne 7 -and $_ -ne 10 -and $_ -ne 14 -and $_ -ne 27 -and $_ -ne 30 }
foreach ($r in $vStyledZeroRows) {
    $cell = $ws.Range("V$r")
    $cell.Value = 0
    $cell.NumberFormat = "0.00"
}

# --- Real (non-zero) values for the provinces that already have FTT data --
# These cells keep the worksheet's default (general) style -- only .Value is
# set, no NumberFormat call, so no extra style gets minted.

# Row 7 - Central Java
$ws.Range("V7").Value = -35.946843853820596

# Row 10 - East Java (U10 keeps the "0.00" style with a real value)
$ws.Range("U10").Value = -65.550286084840735
$ws.Range("U10").NumberFormat = "0.00"
$ws.Range("V10").Value = -37.291246100201867

# Row 14 - Jakarta Special Capital Region
$ws.Range("U14").Value = -66.488825953857457
$ws.Range("V14").Value = -37.738771295818275

# Row 27 - South Sulawesi
$ws.Range("U27").Value = -54.007451575054375
$ws.Range("V27").Value = -5.6249999999999982

# Row 30 - West Java (U30 keeps the "0.00" style with value 0)
$ws.Range("U30").Value = 0
$ws.Range("U30").NumberFormat = "0.00"
$ws.Range("V30").Value = -17.668414683340053

# --- Restore the selection to where the author last left it ---------------
$ws.Range("X34").Select()
